# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer block
# (and the blank paragraph that precedes it) that followed the
# "LOQ4010: Introdução à Engenharia Química (Requisito fraco)" requirement
# line, per the source commit.
$d = $word.ActiveDocument

$marker = "LOQ4010: Introdução à Engenharia Química (Requisito fraco)"

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $marker) {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not locate the '$marker' paragraph"
}

# The three paragraphs to drop are the ones immediately following the
# requirement line: a blank paragraph, "Ver no Jupiter ...", and the
# "(c) 2020 ..." copyright line. Delete them back-to-front so each
# paragraph's Next reference stays valid while the others are removed.
$toDelete = @()
$p = $target.Next()
for ($i = 0; $i -lt 3; $i++) {
    $toDelete += $p
    $p = $p.Next()
}

for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
    $toDelete[$i].Range.Delete()
}
